$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.601.07'
$ws.Range('E2').Value = '  +0.17%  '
$ws.Range('D3').Value = '2.468.86'
$ws.Range('E3').Value = '  -0.61%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '317.28'
$ws.Range('E5').Value = '  +1.21%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '92.20'
$ws.Range('E6').Value = '  -0.32%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.552'
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('E9').Value = '  -0.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '32.92'
$ws.Range('E10').Value = '  +0.10%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0852'
$ws.Range('E11').Value = '  +7.64%  '
$ws.Range('E12').Value = '  +0.31%  '
$ws.Range('D13').Value = '2.848.01'
$ws.Range('E13').Value = '  -0.58%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.89'
$ws.Range('E14').Value = '  -0.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.56'
$ws.Range('E15').Value = '  -4.92%  '
$ws.Range('D16').Value = '2.454.22'
$ws.Range('E16').Value = '  -1.06%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.791'
$ws.Range('E17').Value = '  +1.64%  '
$ws.Range('D18').Value = '41.551.34'
$ws.Range('B19').Value = 'Uniswap'
$ws.Range('C19').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.45'
$ws.Range('E19').Value = '  -1.96%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0948'
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.13'
$ws.Range('E21').Value = '  -2.10%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.31'
$ws.Range('E22').Value = '  +0.86%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '238.93'
$ws.Range('E23').Value = '  +0.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.74'
$ws.Range('E24').Value = '  +0.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.93'
$ws.Range('E25').Value = '  +0.78%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.60'
$ws.Range('E27').Value = '  -1.48%  '
$ws.Range('E28').Value = '  +2.60%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.86'
$ws.Range('E29').Value = '  +1.50%  '
$ws.Range('E30').Value = '  +0.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '161.38'
$ws.Range('E31').Value = '  +2.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.51'
$ws.Range('E32').Value = '  +0.68%  '
$ws.Range('E33').Value = '  +0.09%  '
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('E35').Value = '  +0.74%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '17.22'
$ws.Range('E36').Value = '  -1.01%  '
$ws.Range('E37').Value = '  -0.88%  '
$ws.Range('B38').Value = 'ARBITRUM'
$ws.Range('C38').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.84'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('B39').Value = 'Stellar'
$ws.Range('C39').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.116'
$ws.Range('E39').Value = '  +1.37%  '
$ws.Range('E40').Value = '  -2.83%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.98'
$ws.Range('E41').Value = '  -2.26%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.47'
$ws.Range('E42').Value = '  +3.27%  '
$ws.Range('D43').Value = '1.984.33'
$ws.Range('E43').Value = '  +0.51%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0285'
$ws.Range('E44').Value = '  -0.06%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.86'
$ws.Range('E45').Value = '  -1.55%  '
$ws.Range('E46').Value = '  +0.71%  '
$ws.Range('E47').Value = '  +2.59%  '
$ws.Range('D48').Value = '2.705.32'
$ws.Range('E48').Value = '  -0.59%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '97.21'
$ws.Range('E49').Value = '  -1.17%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '74.11'
$ws.Range('E50').Value = '  +2.17%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '67.24'
$ws.Range('E51').Value = '  -1.68%  '
